$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3059.592
$ws.Range("I15").Value = 3059.592
$ws.Range("K15").Value = 9178.776
$ws.Range("M15").Value = -9009.776
$ws.Range("H80").Value = 12981.889
$ws.Range("I80").Value = 1017
$ws.Range("J80").Value = 16400.428
$ws.Range("K80").Value = 3051
$ws.Range("L80").Value = 49201.284
$ws.Range("M80").Value = -2053
$ws.Range("N80").Value = -51197.284
$ws.Range("H83").Value = 12981.889
$ws.Range("I83").Value = 1017
$ws.Range("J83").Value = 16400.428
$ws.Range("K83").Value = 9153
$ws.Range("L83").Value = 147603.852
$ws.Range("M83").Value = -4161
$ws.Range("N83").Value = -157587.852
$ws.Range("H88").Value = 3500.3076
$ws.Range("J88").Value = 3591.2727
$ws.Range("L88").Value = 3591.2727
$ws.Range("N88").Value = -4403.2727
$ws.Range("H91").Value = 3500.3076
$ws.Range("J91").Value = 3591.2727
$ws.Range("L91").Value = 3591.2727
$ws.Range("N91").Value = -6399.2727
$ws.Range("H138").Value = 2545.476
$ws.Range("I138").Value = 1587.5454
$ws.Range("J138").Value = 3599.2
$ws.Range("K138").Value = 4762.6362
$ws.Range("L138").Value = 10797.6
$ws.Range("M138").Value = 377.3638000000001
$ws.Range("N138").Value = -21077.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2265.9812
$ws.Range("I31").Value = 1917.129
$ws.Range("J31").Value = 2757.5454
$ws.Range("K31").Value = 1917.129
$ws.Range("L31").Value = 2757.5454
$ws.Range("M31").Value = -1622.129
$ws.Range("N31").Value = -3347.5454
$ws.Range("H34").Value = 2265.9812
$ws.Range("I34").Value = 1917.129
$ws.Range("J34").Value = 2757.5454
$ws.Range("K34").Value = 1917.129
$ws.Range("L34").Value = 2757.5454
$ws.Range("M34").Value = -1715.129
$ws.Range("N34").Value = -3161.5454
$ws.Range("H58").Value = 950978.7
$ws.Range("J58").Value = 2785.6
$ws.Range("L58").Value = 2785.6
$ws.Range("N58").Value = -3191.6
$ws.Range("H136").Value = 950978.7
$ws.Range("J136").Value = 2785.6
$ws.Range("L136").Value = 8356.799999999999
$ws.Range("N136").Value = -13456.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 833.6
$ws.Range("J68").Value = 765.8182
$ws.Range("L68").Value = 2297.4546
$ws.Range("N68").Value = -3919.4546
$ws.Range("H69").Value = 1337.5
$ws.Range("J69").Value = 1616.6666
$ws.Range("L69").Value = 4849.9998
$ws.Range("N69").Value = -6471.9998
$ws.Range("H71").Value = 833.6
$ws.Range("J71").Value = 765.8182
$ws.Range("L71").Value = 6892.3638
$ws.Range("N71").Value = -15004.3638
$ws.Range("H72").Value = 1337.5
$ws.Range("J72").Value = 1616.6666
$ws.Range("L72").Value = 14549.9994
$ws.Range("N72").Value = -22661.9994
$ws.Range("H80").Value = 5175.357
$ws.Range("I80").Value = 6362.5
$ws.Range("K80").Value = 19087.5
$ws.Range("M80").Value = -18151.5
$ws.Range("H83").Value = 5175.357
$ws.Range("I83").Value = 6362.5
$ws.Range("K83").Value = 57262.5
$ws.Range("M83").Value = -52582.5
$ws.Range("H129").Value = 3573006
$ws.Range("I129").Value = 521.6
$ws.Range("J129").Value = 5557719.5
$ws.Range("K129").Value = 1564.8
$ws.Range("L129").Value = 16673158.5
$ws.Range("M129").Value = 3435.2
$ws.Range("N129").Value = -16683158.5
$ws.Range("H131").Value = 1009.88
$ws.Range("I131").Value = 562.25
$ws.Range("J131").Value = 1028.5312
$ws.Range("K131").Value = 1686.75
$ws.Range("L131").Value = 3085.5936
$ws.Range("M131").Value = 3353.25
$ws.Range("N131").Value = -13165.5936

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 36200
$ws.Range("J111").Value = 36200
$ws.Range("L111").Value = 36200
$ws.Range("N111").Value = -42334
$ws.Range("H132").Value = 1998.079
$ws.Range("I132").Value = 1293.3846
$ws.Range("J132").Value = 2364.52
$ws.Range("K132").Value = 3880.1538
$ws.Range("L132").Value = 7093.559999999999
$ws.Range("M132").Value = -1350.1538
$ws.Range("N132").Value = -12153.56

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 783.1667
$ws.Range("I22").Value = 832.3333
$ws.Range("J22").Value = 734
$ws.Range("K22").Value = 832.3333
$ws.Range("L22").Value = 734
$ws.Range("M22").Value = -537.3333
$ws.Range("N22").Value = -1324
$ws.Range("H27").Value = 783.1667
$ws.Range("I27").Value = 832.3333
$ws.Range("J27").Value = 734
$ws.Range("K27").Value = 832.3333
$ws.Range("L27").Value = 734
$ws.Range("M27").Value = -725.3333
$ws.Range("N27").Value = -948
$ws.Range("H61").Value = 4681.8184
$ws.Range("I61").Value = 4987.5
$ws.Range("J61").Value = 3866.6667
$ws.Range("K61").Value = 4987.5
$ws.Range("L61").Value = 3866.6667
$ws.Range("M61").Value = -4785.5
$ws.Range("N61").Value = -4270.6667
$ws.Range("H113").Value = 4681.8184
$ws.Range("I113").Value = 4987.5
$ws.Range("J113").Value = 3866.6667
$ws.Range("K113").Value = 4987.5
$ws.Range("L113").Value = 3866.6667
$ws.Range("M113").Value = -2817.5
$ws.Range("N113").Value = -8206.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12875.333
$ws.Range("J45").Value = 12875.333
$ws.Range("L45").Value = 12875.333
$ws.Range("N45").Value = -13857.333
$ws.Range("H46").Value = 60626.348
$ws.Range("J46").Value = 60626.348
$ws.Range("L46").Value = 60626.348
$ws.Range("N46").Value = -61088.348
$ws.Range("H96").Value = 3333.3333
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -6246
$ws.Range("H123").Value = 29234.62
$ws.Range("J123").Value = 29234.62
$ws.Range("L123").Value = 29234.62
$ws.Range("N123").Value = -39034.62
$ws.Range("H132").Value = 1227.341
$ws.Range("I132").Value = 785.17645
$ws.Range("J132").Value = 2730.7
$ws.Range("K132").Value = 2355.52935
$ws.Range("L132").Value = 8192.099999999999
$ws.Range("M132").Value = 174.4706499999998
$ws.Range("N132").Value = -13252.1
$ws.Range("H134").Value = 60626.348
$ws.Range("J134").Value = 60626.348
$ws.Range("L134").Value = 181879.044
$ws.Range("N134").Value = -186949.044
